# BOT; UPDATE DATA
# Appends one new day (2020-05-26, serial 43977) of data to each of the
# three data sheets ("all", "kobe", "other"), matching the upstream
# kansensya_pcr.xlsx refresh.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "all": append a brand-new last row (50) below the existing data
# (row 49 is the current last row; there is no trailing footnote row on
# this sheet, so this is a pure append, not an insert-before).
# ---------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("all")

# Clone the formatting of the row above (48) - rows 44-48 share one
# formatting pattern and the new row continues it - then overwrite with
# this day's figures.
$wsAll.Range("A48:H48").Copy()
$wsAll.Range("A50:H50").PasteSpecial(-4122)
$wsAll.Range("A50").Value = 43977
$wsAll.Range("B50").Value = 285
$wsAll.Range("C50").Value = 282
$wsAll.Range("D50").Value = 20
$wsAll.Range("E50").Value = 17
$wsAll.Range("F50").Value = 3
$wsAll.Range("G50").Value = 12
$wsAll.Range("H50").Value = 250

# ---------------------------------------------------------------------
# Sheet "kobe": insert a new data row just above the trailing footnote
# row (currently row 104) so the footnote shifts down to 105 keeping
# its own formatting/content intact, and the new row inherits the
# formatting of the row that is now above it (103).
# ---------------------------------------------------------------------
$wsKobe = $wb.Worksheets.Item("kobe")
$wsKobe.Rows.Item(104).Insert(-4121)

$wsKobe.Range("A104").Value = 43977
$wsKobe.Range("B104").Value = 25
$wsKobe.Range("C104").Value = 3080
$wsKobe.Range("E104").Value = 285
$wsKobe.Range("F104").Value = 17
$wsKobe.Range("G104").Value = 15
$wsKobe.Range("H104").Value = 2
$wsKobe.Range("I104").Value = 12
$wsKobe.Range("J104").Value = 239

# ---------------------------------------------------------------------
# Sheet "other": same pattern - insert above the footnote row (79),
# which shifts it down to 80, and fill the new row 79 with this day's
# figures.
# ---------------------------------------------------------------------
$wsOther = $wb.Worksheets.Item("other")
$wsOther.Rows.Item(79).Insert(-4121)

$wsOther.Range("A79").Value = 43977
$wsOther.Range("B79").Value = 0
$wsOther.Range("C79").Value = 14
$wsOther.Range("D79").Value = 3
$wsOther.Range("E79").Value = 2
$wsOther.Range("F79").Value = 1
$wsOther.Range("G79").Value = 0
$wsOther.Range("H79").Value = 11

# ---------------------------------------------------------------------
# Restore the per-sheet selections to match the refreshed view (each
# sheet's bottom-right pane now points at a cell near the new last
# row). Do "all" last so it remains the active/selected tab, as it was
# before the edit.
# ---------------------------------------------------------------------
$wsKobe.Range("E108").Select()
$wsOther.Range("C85").Select()
$wsAll.Range("J45").Select()
